# Each adjacent row-pair/trio in the stock report had its Code/Rate/Qty/Value
# columns (B, D, E, F, G) cyclically rotated one row down (wrapping to the
# top of its own group) while the Sl.No (A) and Description (C) stayed put.
# Apply the resulting per-cell values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B142").Value = 63902
$ws.Range("E142").Value = 34.04
$ws.Range("F142").Value = 2
$ws.Range("G142").Value = 64.04000000000001

$ws.Range("B143").Value = 48654
$ws.Range("E143").Value = 38.26
$ws.Range("F143").Value = -1
$ws.Range("G143").Value = -32.02

$ws.Range("B256").Value = 64979
$ws.Range("E256").Value = 314.41
$ws.Range("F256").Value = 82
$ws.Range("G256").Value = 24251.5

$ws.Range("B257").Value = 48719
$ws.Range("E257").Value = 353.35
$ws.Range("F257").Value = -81
$ws.Range("G257").Value = -23955.75

$ws.Range("B271").Value = 48706
$ws.Range("E271").Value = 39.8
$ws.Range("F271").Value = -144
$ws.Range("G271").Value = -4795.2

$ws.Range("B272").Value = 64973
$ws.Range("E272").Value = 35.4
$ws.Range("F272").Value = 150
$ws.Range("G272").Value = 4995

$ws.Range("B308").Value = 57077
$ws.Range("D308").Value = 93.08
$ws.Range("E308").Value = 111.2
$ws.Range("F308").Value = 1
$ws.Range("G308").Value = 93.08

$ws.Range("B309").Value = 61610
$ws.Range("D309").Value = 102.71
$ws.Range("E309").Value = 122.71
$ws.Range("F309").Value = -58
$ws.Range("G309").Value = -5957.18

$ws.Range("B310").Value = 63565
$ws.Range("E310").Value = 109.19
$ws.Range("F310").Value = 60
$ws.Range("G310").Value = 6162.6

$ws.Range("B342").Value = 57802
$ws.Range("E342").Value = 162.71
$ws.Range("F342").Value = -79
$ws.Range("G342").Value = -11334.92

$ws.Range("B343").Value = 63571
$ws.Range("E343").Value = 152.53
$ws.Range("F343").Value = 29
$ws.Range("G343").Value = 4160.92

$ws.Range("B347").Value = 55356
$ws.Range("E347").Value = 54.04
$ws.Range("F347").Value = -158
$ws.Range("G347").Value = -7527.12

$ws.Range("B348").Value = 63510
$ws.Range("E348").Value = 50.66
$ws.Range("F348").Value = 167
$ws.Range("G348").Value = 7955.88

$ws.Range("B367").Value = 63563
$ws.Range("E367").Value = 119.04
$ws.Range("F367").Value = 15
$ws.Range("G367").Value = 1679.4

$ws.Range("B368").Value = 61605
$ws.Range("E368").Value = 133.78
$ws.Range("F368").Value = -13
$ws.Range("G368").Value = -1455.48

$ws.Range("B374").Value = 63560
$ws.Range("E374").Value = 134.87
$ws.Range("F374").Value = 104
$ws.Range("G374").Value = 13193.44

$ws.Range("B375").Value = 60325
$ws.Range("E375").Value = 151.57
$ws.Range("F375").Value = -102
$ws.Range("G375").Value = -12939.72

$ws.Range("B413").Value = 63008
$ws.Range("F413").Value = 504
$ws.Range("G413").Value = 76189.67999999999

$ws.Range("B414").Value = 57857
$ws.Range("F414").Value = 3
$ws.Range("G414").Value = 453.51

$ws.Range("B528").Value = 58047
$ws.Range("D528").Value = 105.54
$ws.Range("E528").Value = 126.1
$ws.Range("F528").Value = 54
$ws.Range("G528").Value = 5699.16

$ws.Range("B529").Value = 47097
$ws.Range("D529").Value = 112.28
$ws.Range("E529").Value = 134.16
$ws.Range("F529").Value = 15
$ws.Range("G529").Value = 1684.2

$ws.Range("B571").Value = 53757
$ws.Range("E571").Value = 16.08
$ws.Range("F571").Value = -159
$ws.Range("G571").Value = -2138.55

$ws.Range("B572").Value = 65069
$ws.Range("E572").Value = 14.3
$ws.Range("F572").Value = 172
$ws.Range("G572").Value = 2313.4

$ws.Range("B591").Value = 45709
$ws.Range("E591").Value = 15.69
$ws.Range("F591").Value = -300
$ws.Range("G591").Value = -3945

$ws.Range("B592").Value = 64925
$ws.Range("E592").Value = 13.97
$ws.Range("F592").Value = 302
$ws.Range("G592").Value = 3971.3

$ws.Range("B596").Value = 53595
$ws.Range("E596").Value = 17.61
$ws.Range("F596").Value = -335
$ws.Range("G596").Value = -4934.55

$ws.Range("B597").Value = 65067
$ws.Range("E597").Value = 15.65
$ws.Range("F597").Value = 338
$ws.Range("G597").Value = 4978.74

$ws.Range("B701").Value = 60025
$ws.Range("E701").Value = 37.22
$ws.Range("F701").Value = -98
$ws.Range("G701").Value = -3217.34

$ws.Range("B702").Value = 64833
$ws.Range("E702").Value = 34.9
$ws.Range("F702").Value = 99
$ws.Range("G702").Value = 3250.17

$ws.Range("B707").Value = 64836
$ws.Range("E707").Value = 104.71
$ws.Range("F707").Value = 7
$ws.Range("G707").Value = 689.5

$ws.Range("B708").Value = 60031
$ws.Range("E708").Value = 111.69
$ws.Range("F708").Value = -5
$ws.Range("G708").Value = -492.5

$ws.Range("B712").Value = 60022
$ws.Range("E712").Value = 37.22
$ws.Range("F712").Value = -113
$ws.Range("G712").Value = -3709.79

$ws.Range("B713").Value = 64830
$ws.Range("E713").Value = 34.9
$ws.Range("F713").Value = 117
$ws.Range("G713").Value = 3841.11

$ws.Range("B864").Value = 54751
$ws.Range("E864").Value = 46.34
$ws.Range("F864").Value = -19
$ws.Range("G864").Value = -776.53

$ws.Range("B865").Value = 65079
$ws.Range("E865").Value = 43.44
$ws.Range("F865").Value = 21
$ws.Range("G865").Value = 858.27
